$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "...le F Beta Score avec Beta fixé suivant certaines hypothèses non
#  confirmées par le métier. "
# becomes two runs:
#   "...le F" / "1 Score."
$search1 = "Nous avons effectué notre modélisation sur la base d’une hypothèse forte : la définition d’une métrique d’évaluation : le F Beta Score avec Beta fixé suivant certaines hypothèses non confirmées par le métier. "
$rng1 = $d.Content
$rng1.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1 = $d.Range($rng1.Start, $rng1.End)
$xml1 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:r><w:t>Nous avons effectué notre modélisation sur la base d’une hypothèse forte : la définition d’une métrique d’évaluation : le F</w:t></w:r>" +
        "<w:r><w:t>1 Score.</w:t></w:r>" +
        "</w:p>"
$r1.InsertXML($xml1)

# --- Change 2 ---------------------------------------------------------
# "L’interprétabilité du modèle pourrait être étoffée en considérant les
#  variables issues du one hot encoding comme une seule et même variable
#  dans la perturbation (un client ne pouvant cumuler plusieurs
#  caractéristiques dans la logique du jeu de données initial. "
# becomes three runs:
#   "L’interprétabilité du modèle pourrait être " / "étoffée et" /
#   " éclaircie (qu’il y a-t-il derrière les sources externes ?)"
$search2 = "L’interprétabilité du modèle pourrait être étoffée en considérant les variables issues du one hot encoding comme une seule et même variable dans la perturbation (un client ne pouvant cumuler plusieurs caractéristiques dans la logique du jeu de données initial. "
$rng2 = $d.Content
$rng2.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2 = $d.Range($rng2.Start, $rng2.End)
$xml2 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:r><w:t xml:space='preserve'>L’interprétabilité du modèle pourrait être </w:t></w:r>" +
        "<w:r><w:t>étoffée et</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> éclaircie (qu’il y a-t-il derrière les sources externes ?)</w:t></w:r>" +
        "</w:p>"
$r2.InsertXML($xml2)
